$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cells whose value must change from 0 to 1
$cellsToSet = @(
    "C5", "E5", "H5", "C7", "E7", "H7", "D8", "H8", "D11", "E11",
    "H11", "H14", "D18", "C26", "E26", "H26", "D29", "H29", "C31", "H31",
    "C32", "E32", "H32", "C33", "E33", "H33", "C34", "E34", "H34", "E35",
    "H35", "D38", "C39", "H39", "C40", "E40", "H40", "C41", "E41", "H41",
    "D51", "H51", "H55", "C56", "E56", "H56", "C58", "H58", "D60", "E61",
    "H61", "D62", "C67", "H67", "C68", "E68", "H68", "E70", "H70", "C72",
    "E72", "H72", "C73", "E73", "H73", "D74", "E74", "H74", "E77", "H77",
    "D86", "H86", "C87", "E87", "H87", "C89", "E89", "H89", "C94", "H94",
    "E96", "H96", "D98", "C99", "E99", "H99", "D100", "H100"
)

foreach ($cellRef in $cellsToSet) {
    $ws.Range($cellRef).Value = 1
}
